$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = -8.403
$ws.Range("D12").Value = -7.290000000000001
$ws.Range("D18").Value = -8.318999999999999
$ws.Range("D37").Value = -8.081999999999999
$ws.Range("D55").Value = -8.319000000000001
$ws.Range("D68").Value = -7.111
$ws.Range("D77").Value = -7.779000000000001
$ws.Range("D78").Value = -8.481
